# Update "想去人数" (want-to-go count) figures in the 北京-漫展信息 workbook
# to reflect the newly scraped data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions) sheet
$wsExhibit.Range("F2").Value  = 332
$wsExhibit.Range("F14").Value = 102
$wsExhibit.Range("F17").Value = 38
$wsExhibit.Range("F18").Value = 1792
$wsExhibit.Range("F24").Value = 1462
$wsExhibit.Range("F28").Value = 626
$wsExhibit.Range("F30").Value = 2478
$wsExhibit.Range("F36").Value = 200
$wsExhibit.Range("F37").Value = 928

# 演出 (Shows) sheet
$wsShow.Range("F23").Value = 110

# 全部类型 (All types) sheet - aggregated view of the same events
$wsAll.Range("F15").Value = 102
$wsAll.Range("F20").Value = 38
$wsAll.Range("F22").Value = 1792
$wsAll.Range("F30").Value = 1462
$wsAll.Range("F35").Value = 626
$wsAll.Range("F41").Value = 928
$wsAll.Range("F44").Value = 110
